$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two embedded logo pictures (and their drawing anchors) ---
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# --- Drop the old rows 2:3 entirely so their custom row-heights (tied to the
#     removed pictures) don't survive into the rebuilt table ---
$ws.Rows("2:3").Delete()

# --- Rebuild the table content. Column A is written top-to-bottom first,
#     then the remaining header/data cells, so that new shared-string
#     entries land in the same order as the target workbook. ---
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "James Bond"
$ws.Range("A3").Value = "Ellen Louise Ripley"
$ws.Range("A4").Value = "Captain Jack Sparrow"

$ws.Range("B4").Value = "31.01.1753"

$ws.Range("C1").Value = "random_int"
$ws.Range("B1").Value = "birthday"

$ws.Range("B2").Value = 45
$ws.Range("C2").Value = 4573

$ws.Range("B3").Value = 102128
$ws.Range("C3").Value = 982630

$ws.Range("C4").Value = 7239

# --- Number formats: birthday column holds real dates (serial numbers)
#     formatted with the built-in short-date format (numFmtId 14) ---
$ws.Range("B2:B4").NumberFormat = "mm-dd-yy"

# --- Header row: bold Calibri font ---
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Bold = $true

# --- Column widths (best-fit) ---
$ws.Columns("A:C").AutoFit()

# --- Selection / active cell ---
$ws.Range("A7").Select() | Out-Null
